# Update the multiplication problems in the table to the new set of values.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "701×7="; New = "288×3=" },
    @{ Old = "826×6="; New = "121×6=" },
    @{ Old = "528×7="; New = "446×6=" },
    @{ Old = "909×2="; New = "553×6=" },
    @{ Old = "710×4="; New = "312×7=" },
    @{ Old = "195×6="; New = "901×3=" },
    @{ Old = "311×5="; New = "148×7=" },
    @{ Old = "140×7="; New = "160×3=" },
    @{ Old = "225×4="; New = "712×5=" },
    @{ Old = "955×2="; New = "274×7=" },
    @{ Old = "155×4="; New = "844×8=" },
    @{ Old = "626×9="; New = "221×2=" },
    @{ Old = "326×6="; New = "304×2=" },
    @{ Old = "199×6="; New = "568×7=" },
    @{ Old = "236×6="; New = "493×3=" },
    @{ Old = "423×4="; New = "488×8=" },
    @{ Old = "345×2="; New = "368×7=" },
    @{ Old = "490×4="; New = "744×7=" },
    @{ Old = "160×5="; New = "955×6=" },
    @{ Old = "349×6="; New = "817×6=" },
    @{ Old = "239×9="; New = "654×9=" },
    @{ Old = "509×8="; New = "825×4=" },
    @{ Old = "691×6="; New = "597×6=" },
    @{ Old = "465×8="; New = "220×4=" },
    @{ Old = "419×5="; New = "114×9=" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
